$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V_Normal")
$ws.Range("F1").Value = "N"
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(V_Normal!`$F`$1,V_Normal!`$A`$2:`$A`$10,V_Normal!`$F`$2:`$F`$10,1)"
Write-Output $s1.Name
Write-Output $s1.Formula
